# Update crypto price/volume figures per the Sat May 20 04:13:50 UTC 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell "D2" "26.907.84"
Set-TextCell "E2" "  +0.22%  "
Set-TextCell "D3" "1.816.16"
Set-TextCell "E3" "  +0.81%  "
Set-TextCell "E4" "  +0.13%  "
Set-TextCell "D5" "309.37"
Set-TextCell "E5" "  +0.07%  "
Set-TextCell "D6" "1.002"
Set-TextCell "E6" "  +0.11%  "
Set-TextCell "D7" "0.4649"
Set-TextCell "E7" "  -0.15%  "
Set-TextCell "D8" "0.3659"
Set-TextCell "E8" "  -0.86%  "
Set-TextCell "D9" "0.07350"
Set-TextCell "E9" "  -0.21%  "
Set-TextCell "D10" "0.8699"
Set-TextCell "E10" "  +0.04%  "
Set-TextCell "D11" "20.31"
Set-TextCell "E11" "  -0.12%  "
Set-TextCell "D12" "1.830.89"
Set-TextCell "E12" "  +2.90%  "
Set-TextCell "D13" "5.381"
Set-TextCell "E13" "  +0.50%  "
Set-TextCell "D14" "0.07094"
Set-TextCell "E14" "  +0.98%  "
Set-TextCell "D15" "6.506"
Set-TextCell "E15" "  +0.13%  "
Set-TextCell "D16" "91.53"
Set-TextCell "E16" "  -0.95%  "
Set-TextCell "D17" "1.003"
Set-TextCell "E17" "  +0.12%  "
Set-TextCell "D18" "0.000008717"
Set-TextCell "E18" "  +0.34%  "
Set-TextCell "E19" "  +0.13%  "
Set-TextCell "D20" "14.65"
Set-TextCell "E20" "  -0.16%  "
Set-TextCell "D21" "26.938.09"
Set-TextCell "D22" "5.298"
Set-TextCell "E22" "  +0.24%  "
Set-TextCell "D23" "10.63"
Set-TextCell "E23" "  +0.31%  "
Set-TextCell "D24" "2.059.30"
Set-TextCell "E24" "  +2.88%  "
Set-TextCell "D25" "1.894"
Set-TextCell "E25" "  -0.75%  "
Set-TextCell "D26" "150.89"
Set-TextCell "E26" "  -0.51%  "
Set-TextCell "D27" "18.29"
Set-TextCell "E27" "  -0.09%  "
Set-TextCell "D28" "2.143"
Set-TextCell "E28" "  +0.70%  "
Set-TextCell "D29" "5.260"
Set-TextCell "E29" "  +0.14%  "
Set-TextCell "D30" "115.06"
Set-TextCell "E30" "  -0.91%  "
Set-TextCell "D31" "0.08901"
Set-TextCell "E31" "  -0.13%  "
Set-TextCell "D32" "0.7556"
Set-TextCell "E32" "  -0.42%  "
Set-TextCell "D33" "1.156"
Set-TextCell "E33" "  +0.66%  "
Set-TextCell "D34" "4.489"
Set-TextCell "E34" "  +0.69%  "
Set-TextCell "E35" "  -0.63%  "
Set-TextCell "E36" "  +0.13%  "
Set-TextCell "E37" "  -1.52%  "
Set-TextCell "D38" "0.05283"
Set-TextCell "E38" "  +0.70%  "
Set-TextCell "D39" "0.01946"
Set-TextCell "E39" "  -0.24%  "
Set-TextCell "D40" "2.976"
Set-TextCell "E40" "  +1.73%  "
Set-TextCell "D41" "7.245"
Set-TextCell "E41" "  +0.42%  "
Set-TextCell "D42" "0.5299"
Set-TextCell "E42" "  +0.14%  "
Set-TextCell "D43" "2.283"
Set-TextCell "E43" "  -3.21%  "
Set-TextCell "E44" "  -0.40%  "
Set-TextCell "D45" "8.427"
Set-TextCell "E45" "  -0.66%  "
Set-TextCell "D46" "0.4869"
Set-TextCell "E46" "  -2.62%  "
Set-TextCell "D47" "10.43"
Set-TextCell "E47" "  +0.85%  "
Set-TextCell "E48" "  +0.15%  "
Set-TextCell "B49" "NEARProtocol"
Set-TextCell "C49" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell "D49" "1.662"
Set-TextCell "E49" "  -0.04%  "
Set-TextCell "B50" "Quant"
Set-TextCell "C50" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell "D50" "103.23"
Set-TextCell "E50" "  -0.90%  "
Set-TextCell "D51" "0.06290"
Set-TextCell "E51" "  +0.12%  "
